# Apply updated optimization results to the workbook
$wb = $excel.ActiveWorkbook

# --- Sheet "optimal_decision_variables": refresh column C (Optimal Value) ---
$ws1 = $wb.Worksheets.Item("optimal_decision_variables")
$ws1.Cells.Item(2,3).Value = 8384.833496447134
$ws1.Cells.Item(3,3).Value = 489.5263038192258
$ws1.Cells.Item(4,3).Value = 9243.078307827231
$ws1.Cells.Item(5,3).Value = 2648.118187534795
$ws1.Cells.Item(6,3).Value = 400.000573735179
$ws1.Cells.Item(7,3).Value = 2962.278461306294
$ws1.Cells.Item(8,3).Value = 383.3384502584042
$ws1.Cells.Item(9,3).Value = 3095.796631495638
$ws1.Cells.Item(10,3).Value = 392.392760441091
$ws1.Cells.Item(11,3).Value = 2846.031302449386
$ws1.Cells.Item(12,3).Value = 398.4983490378368
$ws1.Cells.Item(13,3).Value = 3785.368944021246
$ws1.Cells.Item(14,3).Value = 375.9296629292518
$ws1.Cells.Item(15,3).Value = 1930.868472754827
$ws1.Cells.Item(16,3).Value = 5.860006696573327
$ws1.Cells.Item(17,3).Value = 0.2359663736657835

# --- Sheet "default_influent_quality": append new effluent-quality rows 20-26 ---
$ws2 = $wb.Worksheets.Item("default_influent_quality")
$ws2.Cells.Item(20,1).Value = "BOD"
$ws2.Cells.Item(20,2).Value = 379.8279794331
$ws2.Cells.Item(21,1).Value = "COD"
$ws2.Cells.Item(21,2).Value = 1760.4
$ws2.Cells.Item(22,1).Value = "TKN"
$ws2.Cells.Item(22,2).Value = 154.421
$ws2.Cells.Item(23,1).Value = "TN"
$ws2.Cells.Item(23,2).Value = 174.421
$ws2.Cells.Item(24,1).Value = "TP"
$ws2.Cells.Item(24,2).Value = 166.740350277819
$ws2.Cells.Item(25,1).Value = "TSS"
$ws2.Cells.Item(25,2).Value = 1182.69269083295
$ws2.Cells.Item(26,1).Value = "VSS"
$ws2.Cells.Item(26,2).Value = 771.4574127222

# --- Sheet "optimal_predicted_effluent": refresh column B (Predicted Value) ---
$ws3 = $wb.Worksheets.Item("optimal_predicted_effluent")
$ws3.Cells.Item(2,2).Value = 394.4754148999774
$ws3.Cells.Item(3,2).Value = 6280.757441172927
$ws3.Cells.Item(4,2).Value = 187.0646998043873
$ws3.Cells.Item(5,2).Value = 289.947859739954
$ws3.Cells.Item(6,2).Value = 270.8725110022738
$ws3.Cells.Item(7,2).Value = 4698.768004918164
$ws3.Cells.Item(8,2).Value = 3988.225572803003
$ws3.Cells.Item(9,2).Value = 39.89364043063546
$ws3.Cells.Item(10,2).Value = 600.9788532602144
$ws3.Cells.Item(11,2).Value = 28.99399502661289
$ws3.Cells.Item(12,2).Value = 409.339646332451
$ws3.Cells.Item(13,2).Value = 6276.850016676349
$ws3.Cells.Item(14,2).Value = 187.8796405035193
$ws3.Cells.Item(15,2).Value = 289.947859739954
$ws3.Cells.Item(16,2).Value = 270.8725110022738
$ws3.Cells.Item(17,2).Value = 4588.77891082034
$ws3.Cells.Item(18,2).Value = 3992.319547660497
$ws3.Cells.Item(19,2).Value = 43.65241913653011
$ws3.Cells.Item(20,2).Value = 618.4643355138325
$ws3.Cells.Item(21,2).Value = 30.6684621824118
$ws3.Cells.Item(22,2).Value = 427.7581266571682
$ws3.Cells.Item(23,2).Value = 6335.341081520011
$ws3.Cells.Item(24,2).Value = 189.4178783531866
$ws3.Cells.Item(25,2).Value = 289.9478597399535
$ws3.Cells.Item(26,2).Value = 270.8725110022738
$ws3.Cells.Item(27,2).Value = 4698.768004918181
$ws3.Cells.Item(28,2).Value = 4015.862375824167
$ws3.Cells.Item(29,2).Value = 47.35646695509265
$ws3.Cells.Item(30,2).Value = 667.548642508334
$ws3.Cells.Item(31,2).Value = 31.21403177283925
$ws3.Cells.Item(32,2).Value = 473.1026120541186
$ws3.Cells.Item(33,2).Value = 6539.309515059662
$ws3.Cells.Item(34,2).Value = 188.90839412372
$ws3.Cells.Item(35,2).Value = 291.0537885987609
$ws3.Cells.Item(36,2).Value = 274.766466958643
$ws3.Cells.Item(37,2).Value = 4784.380069592936
$ws3.Cells.Item(38,2).Value = 4081.75634757178
$ws3.Cells.Item(39,2).Value = 49.35120742330808
$ws3.Cells.Item(40,2).Value = 763.4436283550272
$ws3.Cells.Item(41,2).Value = 32.94311110607721
$ws3.Cells.Item(42,2).Value = 437.1465480887408
$ws3.Cells.Item(43,2).Value = 6972.642598167297
$ws3.Cells.Item(44,2).Value = 191.7328993723428
$ws3.Cells.Item(45,2).Value = 301.5071123656577
$ws3.Cells.Item(46,2).Value = 282.5454948157827
$ws3.Cells.Item(47,2).Value = 5466.893162209611
$ws3.Cells.Item(48,2).Value = 4261.205583479632
$ws3.Cells.Item(49,2).Value = 48.74411914095188
$ws3.Cells.Item(50,2).Value = 641.4567194007644
$ws3.Cells.Item(51,2).Value = 33.68583731551557
$ws3.Cells.Item(52,2).Value = 11.3942631474335
$ws3.Cells.Item(53,2).Value = 696.1359776860126
$ws3.Cells.Item(54,2).Value = 472.1379909236755
$ws3.Cells.Item(55,2).Value = 11994.87181490902
$ws3.Cells.Item(56,2).Value = 25.1938313394444
$ws3.Cells.Item(57,2).Value = 313.8282519826552
$ws3.Cells.Item(58,2).Value = 154.3649115539875
$ws3.Cells.Item(59,2).Value = 478.2293589669854
$ws3.Cells.Item(60,2).Value = 157.5830534817553
$ws3.Cells.Item(61,2).Value = 420.2407407250849
$ws3.Cells.Item(62,2).Value = 90.79383362054161
$ws3.Cells.Item(63,2).Value = 10520.30909009808
$ws3.Cells.Item(64,2).Value = 80.87232903298505
$ws3.Cells.Item(65,2).Value = 8550.282458838325
$ws3.Cells.Item(66,2).Value = 0.7645012847116146
$ws3.Cells.Item(67,2).Value = 88.48810766169235
$ws3.Cells.Item(68,2).Value = 10.17443014253985
$ws3.Cells.Item(69,2).Value = 1135.382538259153
$ws3.Cells.Item(70,2).Value = 0.536425234769576
$ws3.Cells.Item(71,2).Value = 57.80096361169348
